$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new row at position 4. This pushes the old row 4 ("Number
# of disability persons" + values) down to row 5, and the old row 5
# (merged "Source: ..." row) down to row 6 - merge range follows.
# ------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ------------------------------------------------------------------
# Row 1: title (merged A1:I1)
# ------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Khelvachauri Municipality"
$r1 = $ws.Range("A1:I1")
$r1.Merge()
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4108
$r1.WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ------------------------------------------------------------------
# Row 2: "(End of year, persons)" subtitle - unchanged text/format,
# just drop the custom row height back to the sheet default.
# ------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ------------------------------------------------------------------
# Row 3: blank label cell + year headers (A3 switches font to Sylfaen)
# ------------------------------------------------------------------
$a3 = $ws.Range("A3")
$a3.Font.Name = "Sylfaen"
$a3.Font.Size = 11

# ------------------------------------------------------------------
# Row 4 (new): "family with disabilities Persons "
# ------------------------------------------------------------------
$a4 = $ws.Range("A4")
$a4.Value = "family with disabilities Persons "
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.Bold = $false
$a4.Borders.LineStyle = -4142
$a4.Borders.Item(8).LineStyle = 1
$a4.Borders.Item(8).Weight = 2
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4108
$a4.WrapText = $true

$ws.Range("B4").Value = 814
$ws.Range("C4").Value = 847
$ws.Range("D4").Value = 854
$ws.Range("E4").Value = 937
$ws.Range("F4").Value = 1010
$ws.Range("G4").Value = 1082
$ws.Range("H4").Value = 1135
$ws.Range("I4").Value = 1192

$row4vals = $ws.Range("B4:I4")
$row4vals.NumberFormat = "#\ ##0"
$row4vals.Font.Name = "Arial"
$row4vals.Font.Size = 10
$row4vals.Borders.LineStyle = -4142
$row4vals.HorizontalAlignment = -4131
$ws.Rows.Item(4).RowHeight = 24.75

# ------------------------------------------------------------------
# Row 5 (shifted): "disabilities Persons "
# ------------------------------------------------------------------
$a5 = $ws.Range("A5")
$a5.Value = "disabilities Persons "
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Borders.LineStyle = -4142
$a5.Borders.Item(9).LineStyle = 1
$a5.Borders.Item(9).Weight = 2
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108
$a5.WrapText = $true

$ws.Range("B5").Value = 957
$ws.Range("C5").Value = 1001
$ws.Range("D5").Value = 1016
$ws.Range("E5").Value = 1110
$ws.Range("F5").Value = 1197
$ws.Range("G5").Value = 1279
$ws.Range("H5").Value = 1336
$ws.Range("I5").Value = 1401

$row5vals = $ws.Range("B5:I5")
$row5vals.NumberFormat = "#\ ##0"
$row5vals.Font.Name = "Arial"
$row5vals.Font.Size = 10
$row5vals.Borders.LineStyle = -4142
$row5vals.HorizontalAlignment = -4131

$i5 = $ws.Range("I5")
$i5.Borders.Item(9).LineStyle = 1
$i5.Borders.Item(9).Weight = 2
$ws.Rows.Item(5).RowHeight = 21

# ------------------------------------------------------------------
# Row 6 (shifted): merged "Source: ..." row (merge already moved with
# the row insert above).
# ------------------------------------------------------------------
$a6 = $ws.Range("A6")
$a6.Font.Name = "Arial"
$a6.Font.Size = 9
$a6.Borders.LineStyle = -4142
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108
$a6.WrapText = $true

$row6rest = $ws.Range("B6:H6")
$row6rest.Font.Name = "Arial"
$row6rest.Font.Size = 9
$row6rest.Borders.LineStyle = -4142
$row6rest.Borders.Item(8).LineStyle = 1
$row6rest.Borders.Item(8).Weight = 2
$row6rest.HorizontalAlignment = -4131
$row6rest.VerticalAlignment = -4108
$row6rest.WrapText = $true
$ws.Rows.Item(6).RowHeight = 27.75

# ------------------------------------------------------------------
# Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.81640625

Write-Output "done"
